$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: fix ExerciceProfessionnel casing/text + column widths ---
$elem = $wb.Worksheets.Item("Elements")

$elem.Range("A6").Value = "DESCNonQualifiant.ExerciceProfessionnel"
$elem.Range("B6").Value = "DESCNonQualifiant.ExerciceProfessionnel"
$elem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$elem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$elem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"

# Target raw column width is 33.640625 characters; the COM ColumnWidth
# setter in this runtime quantizes the stored width to the nearest 1/6,
# so 32.833333333333336 is the input that lands on the closest
# representable width (33.666666666666664) to the target.
$elem.Columns.Item(1).ColumnWidth = 32.833333333333336
$elem.Columns.Item(2).ColumnWidth = 32.833333333333336
